$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Formula = "'39.971.05"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Formula = "'2.216.47"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Formula = "'291.90"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Formula = "'87.67"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Formula = "'30.50"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Formula = "'0.0781"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").Formula = "'50.11"
$ws.Range("E12").Value = "  +5.24%  "
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Formula = "'2.561.20"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Formula = "'13.78"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").Formula = "'2.213.65"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Formula = "'0.732"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Formula = "'39.911.43"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Formula = "'11.13"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Formula = "'65.69"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Formula = "'237.30"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("D28").Formula = "'23.15"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Formula = "'9.25"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Formula = "'2.05"
$ws.Range("E30").Value = "  -6.76%  "
$ws.Range("D31").Formula = "'156.85"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Formula = "'31.96"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Formula = "'0.999"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Formula = "'4.97"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Formula = "'0.0713"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Formula = "'2.96"
$ws.Range("E36").Value = "  +5.50%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Formula = "'1.73"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("D41").Formula = "'15.33"
$ws.Range("E41").Value = "  -3.67%  "
$ws.Range("D42").Formula = "'2.107.28"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Formula = "'3.73"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("D44").Formula = "'0.0270"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D47").Formula = "'1.93"
$ws.Range("E47").Value = "  -7.82%  "
$ws.Range("D48").Formula = "'2.69"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").Formula = "'2.432.82"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Formula = "'1.48"
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").Formula = "'88.62"
$ws.Range("E51").Value = "  -0.83%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Formula = "'17.82"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Formula = "'9.89"
$ws.Range("E46").Value = "  -0.64%  "
